$d = $word.ActiveDocument

# Append a new, truly empty paragraph right after the very last paragraph
# in the document body (before the sectPr). We do this with a
# position-based Range anchored at the end of the story and insert a bare
# paragraph mark, which yields a paragraph with no run at all - matching
# how the existing blank paragraphs in this document look.
$endPos = $d.Content.End
$tailRange = $d.Range($endPos, $endPos)
$tailRange.InsertBefore("`r")

# Add a second new paragraph after that blank one, then fill it with the
# new sentence text.
$blankPara = $d.Paragraphs.Last
$blankPara.Range.InsertParagraphAfter()

$textPara = $d.Paragraphs.Last
$textPara.Range.Text = "first additional commit by develop"
